$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing "Problem Statement" heading gets restyled to Heading 1 + taller row
$ws.Range("A5").Style = "Heading 1"
$ws.Rows.Item(5).RowHeight = 14.4

# New "Internal Discussion" section
$ws.Range("A8").Value = "Internal Discussion"
$ws.Range("A8").Style = "Heading 1"
$ws.Rows.Item(8).RowHeight = 14.4

$ws.Range("B10").Value = "I need to find historical projections that compare to actual performance."
$ws.Range("B11").Value = "Some of my battery work may support this."
$ws.Range("B12").Value = "Compare simple model to ODE"

$ws.Range("B13").Select() | Out-Null

$win = $excel.ActiveWindow
$win.ScrollRow = 4
